{"js": "// Commit: \"cambio al archivo de mauricio por omar\"\n//\n// The document originally contains a single Spanish paragraph:\n//   \"Este es un archivo de prueba\"\n//\n// Omar's edit appends, after that paragraph:\n//   1. a blank paragraph\n//   2. a new paragraph with the text\n//      \"Esto es una prueba de que Omar estuvo aqu\u00ed\"\n//\n// Both new paragraphs keep the same es-ES language formatting used\n// throughout the rest of the document.\n\nconst body = context.document.body;\n\n// Anchor on the existing (last) paragraph in the body.\nconst lastParagraph = body.paragraphs.getLast();\nlastParagraph.load(\"text\");\nawait context.sync();\n\n// 1) New blank paragraph right after the existing text.\nconst blankParagraph = lastParagraph.insertParagraph(\"\", \"After\");\nblankParagraph.font.set({ languageId: \"Spanish (Spain)\" });\n\n// 2) New paragraph carrying Omar's sentence.\nconst omarParagraph = blankParagraph.insertParagraph(\n  \"Esto es una prueba de que Omar estuvo aqu\u00ed\",\n  \"After\"\n);\nomarParagraph.font.set({ languageId: \"Spanish (Spain)\" });\n\nawait context.sync();\n", "ps1": "# Commit: \"cambio al archivo de mauricio por omar\"\n#\n# The document originally contains a single Spanish paragraph:\n#   \"Este es un archivo de prueba\"\n#\n# Omar's edit appends, after that paragraph:\n#   1. a blank paragraph\n#   2. a new paragraph with the text\n#      \"Esto es una prueba de que Omar estuvo aqu\u00ed\"\n#\n# Both new paragraphs keep the same es-ES language formatting used\n# throughout the rest of the document (InsertParagraphAfter() carries\n# the formatting of the paragraph it is attached to).\n\n$d = $word.ActiveDocument\n\n# Collapse to the very end of the document, right after the existing text.\n$endRange = $d.Paragraphs.Last.Range\n$endRange.Collapse(0)  # wdCollapseEnd\n\n# 1) New blank paragraph right after the existing text.\n$endRange.InsertParagraphAfter()\n\n# 2) New paragraph carrying Omar's sentence.\n$d.Paragraphs.Last.Range.Collapse(0)\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"Esto es una prueba de que Omar estuvo aqu\u00ed\"\n"}
